$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Remove the "reviews_count" column (column E) entirely, shifting all
# subsequent columns (reviews_average, latitude, longitude,
# is_permanently_closed, gmaps_link, latest_review_date) one position left.
$ws.Range("E1").EntireColumn.Delete()
